$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new column ("ABC") is inserted immediately before the existing "Spcl Allowance"
# column (O), shifting O..X right by one to P..Y.
$ws.Columns("O:O").Insert()

# Match the width of the new column to its left neighbour (N / "HRA"), which is
# what Excel does by default when a whole column is inserted.
$ws.Columns("O:O").ColumnWidth = 5.642857142857143

# Header + data for the newly inserted "ABC" column.
$ws.Range("O1").Value = "ABC"
$ws.Range("O2").Value = 0
